# Updated capital structure database
# - Row 2 (formerly unnamed company "1") becomes company "2" with refreshed
#   capital-structure metrics; roe/roe_cost_equity (W/Y) no longer populated.
# - Row 3 (Pacific Edge Limited) gets refreshed capital-structure metrics.
# - Row 4 is a new entry for "Aroa Biosurgery Limited (ASX:ARX)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Row 2: update existing figures -----
# B2 holds a numeric-looking placeholder label that must stay text (not become
# the number 2); force text format, assign, then drop back to Normal style so
# no stray formatting survives on the cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Value = 0.0722
$ws.Range("G2").Value = -0.09756232686980602
$ws.Range("H2").Value = -0.6388365650969529
$ws.Range("I2").Value = -0.9399612761757091
$ws.Range("J2").Value = -0.9399612761757091
$ws.Range("K2").Value = -23.3
$ws.Range("L2").Value = -1.290858725761773
$ws.Range("U2").Value = 16.55
$ws.Range("V2").Value = 0.01830144863430278
$ws.Range("W2").ClearContents()
$ws.Range("X2").Value = 0.05576330758190857
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").Value = 7.725212935211585
$ws.Range("AA2").Value = -177.0954031034483
$ws.Range("AB2").Value = 0.05502162670965156
$ws.Range("AC2").Value = -177.1504247301579
$ws.Range("AD2").Value = 13.8
$ws.Range("AE2").Value = 0.01650517485776423
$ws.Range("AF2").Value = 13.81650517485776
$ws.Range("AG2").Value = -2.733494825142236
$ws.Range("AH2").Value = 0.01504874936566615
$ws.Range("AI2").Value = 0.1839343449578147
$ws.Range("AJ2").Value = -0.003031939196334803
$ws.Range("AK2").Value = -0.04667334711164749
$ws.Range("AL2").Value = 1.71
$ws.Range("AM2").Value = 1.485
$ws.Range("AN2").Value = -0.9071189114573063
$ws.Range("AO2").Value = -9.929824561403509
$ws.Range("AP2").Value = 0.1796815108882032
$ws.Range("AQ2").Value = -11.43434343434344

# ----- Row 3: update existing figures (Pacific Edge Limited) -----
$ws.Range("D3").Value = 0.0722
$ws.Range("G3").Value = -1.197590361445783
$ws.Range("H3").Value = -2.698795180722891
$ws.Range("I3").Value = -2.698795180722891
$ws.Range("J3").Value = -2.698795180722891
$ws.Range("K3").Value = -11
$ws.Range("L3").Value = -2.650602409638554
$ws.Range("U3").Value = 4.15
$ws.Range("V3").Value = 0.006509803921568628
$ws.Range("W3").Value = -3.559870550161813
$ws.Range("X3").Value = 0.05516479363085694
$ws.Range("Y3").Value = -3.61503534379267
$ws.Range("Z3").Value = 1.788793103448276
$ws.Range("AA3").Value = -4.827586206896552
$ws.Range("AB3").Value = 0.05504841870654482
$ws.Range("AC3").Value = -4.882634625603098
$ws.Range("AD3").Value = 2.9
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 2.9
$ws.Range("AG3").Value = -1.25
$ws.Range("AH3").Value = 0.00452841973766396
$ws.Range("AI3").Value = 0.1260869565217391
$ws.Range("AJ3").Value = -0.001964636542239686
$ws.Range("AK3").Value = -0.06631299734748013
$ws.Range("AM3").Value = -0.203
$ws.Range("AN3").Value = -0.2636363636363636
$ws.Range("AP3").Value = 0.1136363636363637
$ws.Range("AQ3").Value = 55.17241379310344

# ----- Row 4: new entry - Aroa Biosurgery Limited (ASX:ARX) -----
$ws.Range("A4").Value = "New Zealand"
$ws.Range("B4").Value = "Aroa Biosurgery Limited (ASX:ARX)"
$ws.Range("C4").Value = "Drugs (Biotechnology)"
$ws.Range("G4").Value = 0.230863309352518
$ws.Range("H4").Value = -0.02381294964028777
$ws.Range("I4").Value = -0.4148418010770901
$ws.Range("J4").Value = -0.4148418010770901
$ws.Range("K4").Value = -12.3
$ws.Range("L4").Value = -0.8848920863309353
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 12.4
$ws.Range("V4").Value = 0.04647676161919041
$ws.Range("X4").Value = 0.0563618215329602
$ws.Range("Z4").Value = 842.1601176470587
$ws.Range("AA4").Value = -349.36322
$ws.Range("AB4").Value = 0.05499483471275831
$ws.Range("AC4").Value = -349.4182148347127
$ws.Range("AD4").Value = 10.9
$ws.Range("AE4").Value = 0.01650517485776423
$ws.Range("AF4").Value = 10.91650517485776
$ws.Range("AG4").Value = -1.483494825142236
$ws.Range("AH4").Value = 0.0393080892616895
$ws.Range("AI4").Value = 0.2094634921937196
$ws.Range("AJ4").Value = -0.005591415521490204
$ws.Range("AK4").Value = -0.03735209879648099
$ws.Range("AL4").Value = 1.71
$ws.Range("AM4").Value = 1.688
$ws.Range("AN4").Value = -2.587230002373606
$ws.Range("AO4").Value = -3.380116959064328
$ws.Range("AP4").Value = 0.3521231486214659
$ws.Range("AQ4").Value = -3.424170616113744
